# replicate_decision.xlsx - "Just loads of stuff I've forgotten to commit"
#
# Fills in three previously-blank T3 results (C3, C5, C7 go from 0 -> 1)
# on Sheet1, then leaves the selection on C8 (where the user's cursor
# ended up after typing the values down the column).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = 1
$ws.Range("C5").Value = 1
$ws.Range("C7").Value = 1

$ws.Range("C8").Select() | Out-Null
